$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Term Type" header in column O, matching the existing header look
# (bold border, no fill) but with a smaller font and wrapped text.
$ws.Range("A1").Copy()
$ws.Range("O1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("O1").Value = "Term Type"
$ws.Range("O1").Font.Size = 11
$ws.Range("O1").WrapText = $true

$ws.Range("H2:H4").Clear()

$ws.Range("O1").Select() | Out-Null
